# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.329.60"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "2.174.94"
$ws.Range("E3").Value = "  -2.11%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.62%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.67%  "

$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.73%  "

$ws.Range("D14").Value = "2.502.70"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").Value = "2.181.45"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.767"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.47%  "

$ws.Range("D18").Value = "42.238.87"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("E19").Value = "  -3.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.02%  "

$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.39"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.91%  "

$ws.Range("E29").Value = "  -2.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0808"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.00%  "

$ws.Range("E35").Value = "  -1.78%  "

$ws.Range("E36").Value = "  -2.90%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0333"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.59%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.91%  "

$ws.Range("E41").Value = "  -1.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "58.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("E43").Value = "  -7.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.83"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0968"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.36%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.18%  "

$ws.Range("B48").Value = "WOONetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.456"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.31%  "

$ws.Range("E50").Value = "  -1.79%  "

$ws.Range("E51").Value = "  +0.18%  "
